# Auto-generated Excel COM-interop script to apply scheduled market-price refresh
# updates to the Leviathan_Profits workbook (Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 864.43475
$ws.Range("J28").Value = 3373.75
$ws.Range("L28").Value = 3373.75
$ws.Range("N28").Value = -4343.75
$ws.Range("H41").Value = 918.8570999999999
$ws.Range("I41").Value = 996
$ws.Range("J41").Value = 793.5
$ws.Range("K41").Value = 996
$ws.Range("L41").Value = 793.5
$ws.Range("M41").Value = -556
$ws.Range("N41").Value = -1673.5
$ws.Range("H62").Value = 4122.8
$ws.Range("I62").Value = 3714.3572
$ws.Range("J62").Value = 5075.8335
$ws.Range("K62").Value = 3714.3572
$ws.Range("L62").Value = 5075.8335
$ws.Range("M62").Value = -3090.3572
$ws.Range("N62").Value = -6323.8335
$ws.Range("H65").Value = 4122.8
$ws.Range("I65").Value = 3714.3572
$ws.Range("J65").Value = 5075.8335
$ws.Range("K65").Value = 18571.786
$ws.Range("L65").Value = 25379.1675
$ws.Range("M65").Value = -15451.786
$ws.Range("N65").Value = -31619.1675
$ws.Range("H101").Value = 62502280
$ws.Range("I101").Value = 2605.1428
$ws.Range("J101").Value = 500000000
$ws.Range("K101").Value = 7815.428400000001
$ws.Range("L101").Value = 1500000000
$ws.Range("M101").Value = -6193.428400000001
$ws.Range("N101").Value = -1500003244
$ws.Range("H116").Value = 4570.143
$ws.Range("I116").Value = 3995.5
$ws.Range("K116").Value = 3995.5
$ws.Range("M116").Value = -553.5
$ws.Range("H132").Value = 1753.4166
$ws.Range("I132").Value = 1533.0294
$ws.Range("K132").Value = 4599.0882
$ws.Range("M132").Value = -2069.0882

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2442.0908
$ws.Range("I2").Value = 1874.2222
$ws.Range("K2").Value = 1874.2222
$ws.Range("M2").Value = -1761.2222
$ws.Range("H6").Value = 1500
$ws.Range("I6").Value = 1500
$ws.Range("K6").Value = 1500
$ws.Range("M6").Value = -1327
$ws.Range("H45").Value = 6729.4
$ws.Range("I45").Value = 11230.637
$ws.Range("J45").Value = 3192.7144
$ws.Range("K45").Value = 11230.637
$ws.Range("L45").Value = 3192.7144
$ws.Range("M45").Value = -10853.637
$ws.Range("N45").Value = -3946.7144
$ws.Range("H116").Value = 2442.0908
$ws.Range("I116").Value = 1874.2222
$ws.Range("K116").Value = 1874.2222
$ws.Range("M116").Value = 419.7778000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2442.0908
$ws.Range("I3").Value = 1874.2222
$ws.Range("K3").Value = 1874.2222
$ws.Range("M3").Value = -1760.2222
$ws.Range("H62").Value = 27590.5
$ws.Range("J62").Value = 27590.5
$ws.Range("L62").Value = 27590.5
$ws.Range("N62").Value = -28962.5
$ws.Range("H65").Value = 27590.5
$ws.Range("J65").Value = 27590.5
$ws.Range("L65").Value = 82771.5
$ws.Range("N65").Value = -89635.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1146.5
$ws.Range("I22").Value = 1374.75
$ws.Range("J22").Value = 994.3333
$ws.Range("K22").Value = 1374.75
$ws.Range("L22").Value = 994.3333
$ws.Range("M22").Value = -1024.75
$ws.Range("N22").Value = -1694.3333
$ws.Range("H105").Value = 2106
$ws.Range("J105").Value = 1999.5
$ws.Range("L105").Value = 1999.5
$ws.Range("N105").Value = -5493.5
$ws.Range("H122").Value = 113971.336
$ws.Range("I122").Value = 113971.336
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 341914.008
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -339464.008
$ws.Range("N122").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 14483.333
$ws.Range("I3").Value = 14483.333
$ws.Range("K3").Value = 43449.999
$ws.Range("M3").Value = -43337.999
$ws.Range("H108").Value = 960.8570999999999
$ws.Range("I108").Value = 960.8570999999999
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 2882.5713
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -2.57129999999961
$ws.Range("N108").ClearContents()
$ws.Range("H117").Value = 4553.524
$ws.Range("J117").Value = 5030.8237
$ws.Range("L117").Value = 15092.4711
$ws.Range("N117").Value = -21976.4711
$ws.Range("H132").Value = 1952
$ws.Range("J132").Value = 2322.4443
$ws.Range("L132").Value = 20901.9987
$ws.Range("N132").Value = -25961.9987
$ws.Range("H139").Value = 6625.9
$ws.Range("I139").Value = 8343
$ws.Range("K139").Value = 25029
$ws.Range("M139").Value = -19889

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 23570.5
$ws.Range("J52").Value = 25284.6
$ws.Range("L52").Value = 25284.6
$ws.Range("N52").Value = -25802.6
$ws.Range("H53").Value = 30043
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 30043
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 30043
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -31305
$ws.Range("H102").Value = 1573.3448
$ws.Range("I102").Value = 1344.56
$ws.Range("J102").Value = 3003.25
$ws.Range("K102").Value = 1344.56
$ws.Range("L102").Value = 3003.25
$ws.Range("M102").Value = 277.4400000000001
$ws.Range("N102").Value = -6247.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 55999.8
$ws.Range("J7").Value = 5681.8184
$ws.Range("K7").Value = 55999.8
$ws.Range("L7").Value = 5681.8184
$ws.Range("M7").Value = -55887.8
$ws.Range("N7").Value = -5905.8184
$ws.Range("H22").Value = 4432.6
$ws.Range("I22").Value = 8225
$ws.Range("J22").Value = 3484.5
$ws.Range("K22").Value = 8225
$ws.Range("L22").Value = 3484.5
$ws.Range("M22").Value = -7930
$ws.Range("N22").Value = -4074.5
$ws.Range("H27").Value = 4432.6
$ws.Range("I27").Value = 8225
$ws.Range("J27").Value = 3484.5
$ws.Range("K27").Value = 8225
$ws.Range("L27").Value = 3484.5
$ws.Range("M27").Value = -8118
$ws.Range("N27").Value = -3698.5
$ws.Range("H46").Value = 29003.588
$ws.Range("I46").Value = 142984.67
$ws.Range("J46").Value = 4579.0713
$ws.Range("K46").Value = 142984.67
$ws.Range("L46").Value = 4579.0713
$ws.Range("M46").Value = -142796.67
$ws.Range("N46").Value = -4955.0713
$ws.Range("H55").Value = 177.13043
$ws.Range("I55").Value = 164.61539
$ws.Range("K55").Value = 164.61539
$ws.Range("M55").Value = 8.384610000000009
$ws.Range("I126").Value = 55999.8
$ws.Range("J126").Value = 5681.8184
$ws.Range("K126").Value = 167999.4
$ws.Range("L126").Value = 17045.4552
$ws.Range("M126").Value = -165529.4
$ws.Range("N126").Value = -21985.4552

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 45456360
$ws.Range("I107").Value = 2658
$ws.Range("K107").Value = 7974
$ws.Range("M107").Value = -6054
$ws.Range("H131").Value = 92326.25
$ws.Range("J131").Value = 92326.25
$ws.Range("L131").Value = 92326.25
$ws.Range("N131").Value = -102406.25

Write-Host "Applied scheduled market data refresh to Leviathan_Profits workbook."
